# weight measurement exp2 02-06
# Add a new "t28" weight-measurement column (M) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("M1").Value = "t28"

# Match the font size (14pt) used by the other measurement columns (J, K, L)
# for the new column's data cells (M3 is left at default format, matching
# the source data).
$ws.Range("M2,M4:M33").Font.Size = 14

$ws.Range("M2").Value = 17.2
$ws.Range("M3").Value = 18
$ws.Range("M4").Value = 19.8
$ws.Range("M5").Value = 18.899999999999999
$ws.Range("M6").Value = 18.2
$ws.Range("M7").Value = 17.399999999999999
$ws.Range("M8").Value = 17.7
$ws.Range("M9").Value = 18.399999999999999
$ws.Range("M10").Value = 18.2
$ws.Range("M11").Value = 16.3
$ws.Range("M12").Value = 18.3
$ws.Range("M13").Value = 19.2
$ws.Range("M14").Value = 19.100000000000001
$ws.Range("M15").Value = 17.7
$ws.Range("M16").Value = 19.100000000000001
$ws.Range("M17").Value = 18.3
$ws.Range("M18").Value = 21.2
$ws.Range("M19").Value = 15.9
$ws.Range("M20").Value = 18.3
$ws.Range("M21").Value = 17.8
$ws.Range("M22").Value = 17.100000000000001
$ws.Range("M23").Value = 18.399999999999999
$ws.Range("M24").Value = 17.5
$ws.Range("M25").Value = 17.8
$ws.Range("M26").Value = 18.899999999999999
$ws.Range("M27").Value = 17.7
$ws.Range("M28").Value = 18
$ws.Range("M29").Value = 18.2
$ws.Range("M30").Value = 19.3
$ws.Range("M31").Value = 17.7
$ws.Range("M32").Value = 16.600000000000001
$ws.Range("M33").Value = 18.8

# Scroll back to the top and select the newly entered cell, matching the
# author's final view state.
$ws.Range("M4").Select()
